$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.907049298286438
$ws.Range("B1").Value = 1.994985222816467
$ws.Range("C1").Value = 2.021505117416382
$ws.Range("D1").Value = 2.644572257995605
$ws.Range("E1").Value = 3.3934326171875
